# Quarterly database update: shift all quarterly columns (D..M) one quarter to
# the left (dropping the oldest quarter, which fell out of the shared-string
# table) and populate the newly-opened rightmost column (M) with the freshly
# reported quarter's figures. Also nudges the "year-end" wide column marker
# (width 31) one position to the left so it keeps tracking the Q4 column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row 8: "دوره مالی" (financial period) labels for columns D..M
# ---------------------------------------------------------------------------
$ws.Range("D8").Value = "فصل دوم منتهی به 1399/09"
$ws.Range("E8").Value = "فصل سوم منتهی به 1399/12"
$ws.Range("F8").Value = "فصل چهارم منتهی به 1400/03"
$ws.Range("G8").Value = "فصل اول منتهی به 1400/06"
$ws.Range("H8").Value = "فصل دوم منتهی به 1400/09"
$ws.Range("I8").Value = "فصل سوم منتهی به 1400/12"
$ws.Range("J8").Value = "فصل چهارم منتهی به 1401/03"
$ws.Range("K8").Value = "فصل اول منتهی به 1401/06"
$ws.Range("L8").Value = "فصل دوم منتهی به 1401/09"
$ws.Range("M8").Value = "فصل سوم منتهی به 1401/12"

# ---------------------------------------------------------------------------
# Header row 9: "تاریخ انتشار" (publish date) labels for columns D..M
# ---------------------------------------------------------------------------
$ws.Range("D9").Value = "1400-11-27 (4)"
$ws.Range("E9").Value = "1401-01-30 (2)"
$ws.Range("F9").Value = "1401-07-20 (9)"
$ws.Range("G9").Value = "1401-07-30 (2)"
$ws.Range("H9").Value = "1401-12-03 (4)"
$ws.Range("I9").Value = "1402-01-30 (2)"
$ws.Range("J9").Value = "1402-01-30 (7)"
$ws.Range("K9").Value = "1401-07-30"
$ws.Range("L9").Value = "1401-12-03 (2)"
$ws.Range("M9").Value = "1402-01-30"

# ---------------------------------------------------------------------------
# Row 11: فروش (Sales)
# ---------------------------------------------------------------------------
$ws.Range("D11").Value = 3444
$ws.Range("E11").Value = 4732
$ws.Range("F11").Value = 5940
$ws.Range("G11").Value = 6516
$ws.Range("H11").Value = 5751
$ws.Range("I11").Value = 6957
$ws.Range("J11").Value = 6476
$ws.Range("K11").Value = 7554
$ws.Range("L11").Value = 6081
$ws.Range("M11").Value = 8263

# ---------------------------------------------------------------------------
# Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold)
# ---------------------------------------------------------------------------
$ws.Range("D12").Value = -1453
$ws.Range("E12").Value = -2082
$ws.Range("F12").Value = -2280
$ws.Range("G12").Value = -2421
$ws.Range("H12").Value = -2087
$ws.Range("I12").Value = -2739
$ws.Range("J12").Value = -2764
$ws.Range("K12").Value = -3419
$ws.Range("L12").Value = -3365
$ws.Range("M12").Value = -4652

# ---------------------------------------------------------------------------
# Row 13: سود (زیان) ناخالص (Gross profit)
# ---------------------------------------------------------------------------
$ws.Range("D13").Value = 1991
$ws.Range("E13").Value = 2651
$ws.Range("F13").Value = 3661
$ws.Range("G13").Value = 4095
$ws.Range("H13").Value = 3664
$ws.Range("I13").Value = 4217
$ws.Range("J13").Value = 3712
$ws.Range("K13").Value = 4135
$ws.Range("L13").Value = 2717
$ws.Range("M13").Value = 3612

# ---------------------------------------------------------------------------
# Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses)
# ---------------------------------------------------------------------------
$ws.Range("D14").Value = -163
$ws.Range("E14").Value = -140
$ws.Range("F14").Value = -133
$ws.Range("G14").Value = -134
$ws.Range("H14").Value = -240
$ws.Range("I14").Value = -193
$ws.Range("J14").Value = -222
$ws.Range("K14").Value = -275
$ws.Range("L14").Value = -348
$ws.Range("M14").Value = -275

# Row 15: هزینه کاهش ارزش دریافتنی‌ها (always "-") is unchanged.

# ---------------------------------------------------------------------------
# Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی (Other operating income/exp.)
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = 119
$ws.Range("E16").Value = -45
$ws.Range("F16").Value = -31
$ws.Range("G16").Value = "-"
$ws.Range("H16").Value = 691
$ws.Range("I16").Value = -637
$ws.Range("J16").Value = -2
$ws.Range("K16").Value = "-"
$ws.Range("L16").Value = 662
$ws.Range("M16").Value = -380

# ---------------------------------------------------------------------------
# Row 17: سود (زیان) عملیاتی (Operating profit)
# ---------------------------------------------------------------------------
$ws.Range("D17").Value = 1948
$ws.Range("E17").Value = 2466
$ws.Range("F17").Value = 3497
$ws.Range("G17").Value = 3961
$ws.Range("H17").Value = 4115
$ws.Range("I17").Value = 3388
$ws.Range("J17").Value = 3487
$ws.Range("K17").Value = 3860
$ws.Range("L17").Value = 3031
$ws.Range("M17").Value = 2957

# ---------------------------------------------------------------------------
# Row 18: هزینه های مالی (Financial expenses)
# ---------------------------------------------------------------------------
$ws.Range("D18").Value = -9
$ws.Range("E18").Value = -13
$ws.Range("F18").Value = -15
$ws.Range("G18").Value = -214
$ws.Range("H18").Value = -198
$ws.Range("I18").Value = -210
$ws.Range("J18").Value = -266
$ws.Range("K18").Value = -558
$ws.Range("L18").Value = -683
$ws.Range("M18").Value = -811

# ---------------------------------------------------------------------------
# Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی (Other non-op. income/exp.)
# ---------------------------------------------------------------------------
$ws.Range("D19").Value = 2114
$ws.Range("E19").Value = -774
$ws.Range("F19").Value = -180
$ws.Range("G19").Value = 252
$ws.Range("H19").Value = 436
$ws.Range("I19").Value = 576
$ws.Range("J19").Value = 359
$ws.Range("K19").Value = 310
$ws.Range("L19").Value = 1815
$ws.Range("M19").Value = 1063

# ---------------------------------------------------------------------------
# Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات (Pre-tax profit)
# ---------------------------------------------------------------------------
$ws.Range("D20").Value = 4053
$ws.Range("E20").Value = 1680
$ws.Range("F20").Value = 3301
$ws.Range("G20").Value = 3999
$ws.Range("H20").Value = 4353
$ws.Range("I20").Value = 3754
$ws.Range("J20").Value = 3580
$ws.Range("K20").Value = 3611
$ws.Range("L20").Value = 4162
$ws.Range("M20").Value = 3209

# ---------------------------------------------------------------------------
# Row 21: مالیات (Tax)
# ---------------------------------------------------------------------------
$ws.Range("D21").Value = -708
$ws.Range("E21").Value = -326
$ws.Range("F21").Value = 103
$ws.Range("G21").Value = -726
$ws.Range("H21").Value = -393
$ws.Range("I21").Value = -537
$ws.Range("J21").Value = -646
$ws.Range("K21").Value = -630
$ws.Range("L21").Value = 224
$ws.Range("M21").Value = -607

# ---------------------------------------------------------------------------
# Row 22: سود (زیان) خالص عملیات در حال تداوم (Net profit from continuing ops)
# ---------------------------------------------------------------------------
$ws.Range("D22").Value = 3345
$ws.Range("E22").Value = 1353
$ws.Range("F22").Value = 3404
$ws.Range("G22").Value = 3273
$ws.Range("H22").Value = 3960
$ws.Range("I22").Value = 3217
$ws.Range("J22").Value = 2934
$ws.Range("K22").Value = 2982
$ws.Range("L22").Value = 4386
$ws.Range("M22").Value = 2602

# Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی (always "-") unchanged.

# ---------------------------------------------------------------------------
# Row 24: سود (زیان) خالص (Net profit) - mirrors row 22
# ---------------------------------------------------------------------------
$ws.Range("D24").Value = 3345
$ws.Range("E24").Value = 1353
$ws.Range("F24").Value = 3404
$ws.Range("G24").Value = 3273
$ws.Range("H24").Value = 3960
$ws.Range("I24").Value = 3217
$ws.Range("J24").Value = 2934
$ws.Range("K24").Value = 2982
$ws.Range("L24").Value = 4386
$ws.Range("M24").Value = 2602

# Row 25: سود هر سهم پس از کسر مالیات (always 0) unchanged.

# ---------------------------------------------------------------------------
# Row 26: سرمایه (Capital)
# ---------------------------------------------------------------------------
$ws.Range("D26").Value = 10598
$ws.Range("E26").Value = 11792
$ws.Range("F26").Value = 12345
$ws.Range("G26").Value = 11034
$ws.Range("H26").Value = 10115
$ws.Range("I26").Value = 10456
$ws.Range("J26").Value = 9801
$ws.Range("K26").Value = 9286
$ws.Range("L26").Value = 8289
$ws.Range("M26").Value = 16501

# Row 27: سود هر سهم بر اساس آخرین سرمایه (always 0) unchanged.

# ---------------------------------------------------------------------------
# Column widths: the "year-end" (فصل چهارم) quarter column is drawn wider
# (31) than the rest (29). Since the data shifted one quarter to the left,
# that wide marker moves from columns G/K to columns F/J.
# ---------------------------------------------------------------------------
$narrowWidth = $ws.Columns.Item(4).ColumnWidth
$wideWidth = $ws.Columns.Item(7).ColumnWidth

$ws.Columns.Item(6).ColumnWidth = $wideWidth
$ws.Columns.Item(7).ColumnWidth = $narrowWidth
$ws.Columns.Item(10).ColumnWidth = $wideWidth
$ws.Columns.Item(11).ColumnWidth = $narrowWidth
